$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: sex value for PQY was unparseable ("F" couldn't be classified) -> captured as numeric 2
$ws.Range("D3").Value = 2

# New row 6: an additional captured record (added first so its new shared
# string "prakhar yadav" lands before "2Six" in the shared-strings table)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "prakhar yadav"
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = "M"

# Row 5: age value for LS was unparseable ("2Six") -> captured as the literal string
$ws.Range("C5").Value = "2Six"

# Update the active selection to reflect where the edit was made
$ws.Range("C5").Select()
